# chore: simulator full-month coverage, persist logs, fix employees
#
# - Replaces the placeholder client names with the real client names on
#   both the "Weekly Timesheet" and "Jason Schema" sheets.
# - Fixes Phil Henderson's employee id on the "Jason Schema" sheet.
# - Fills in the simulator output (rate/total columns) that was previously
#   stubbed out at 0 now that the simulator covers the full month.

$wb = $excel.ActiveWorkbook

$timesheet = $wb.Worksheets.Item("Weekly Timesheet")
$schema = $wb.Worksheets.Item("Jason Schema")

# --- Client name corrections (same 5 clients on both sheets) -----------
$clients = @("Winn", "Keevil", "Howard", "Markfield", "Layne")

for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = 2 + $i
    $timesheet.Range("B$row").Value = $clients[$i]
    $schema.Range("D$row").Value = $clients[$i]
}

# --- Employee id correction (Jason Schema, rows 2-6) --------------------
for ($row = 2; $row -le 6; $row++) {
    $schema.Range("B$row").Value = "emp_75yd72zj"
}

# --- Simulator hours now produce real rate/total figures ---------------
# Weekly Timesheet: rows 2-6, Rate (E) and Total (F)
for ($row = 2; $row -le 6; $row++) {
    $timesheet.Range("E$row").Value = 90
    $timesheet.Range("F$row").Value = 720
}

# Weekly Timesheet subtotal / grand total rows
$timesheet.Range("F8").Value = 3600
$timesheet.Range("F11").Value = 3600
$timesheet.Range("F13").Value = 3600

# Jason Schema: rows 2-6, Rate (F) and Total (G)
for ($row = 2; $row -le 6; $row++) {
    $schema.Range("F$row").Value = 90
    $schema.Range("G$row").Value = 720
}
